# Add 2022-Q1 data
#
# Plan:
#  - The current "总计" sheet becomes "2022-Q1" (reusing its sheetId/rId slot),
#    and is repopulated with the fund-holding detail for 2022-Q1.
#  - A brand new sheet named "总计" is inserted right after it, repopulated
#    with the historical roll-up table plus the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force the cell to be stored as text (even when the text looks numeric,
    # e.g. "011429" or "2.83"), then strip the leftover number-format style
    # so the cell keeps the same (unstyled) look as its neighbours.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Turn the existing "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "2022-Q1"
$q1 = $wb.Worksheets.Item("2022-Q1")

$q1.Cells.Clear()

# Header row (column B..H)
Set-TextValue $q1.Range("B1") "基金代码"
Set-TextValue $q1.Range("C1") "基金名称"
Set-TextValue $q1.Range("D1") "基金规模"
Set-TextValue $q1.Range("E1") "股票总仓位"
Set-TextValue $q1.Range("F1") "仓位占比"
Set-TextValue $q1.Range("G1") "持有市值(亿元)"
Set-TextValue $q1.Range("H1") "仓位排名"

# Data rows
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "011429"
Set-TextValue $q1.Range("C2") "前海开源民裕进取混合"
Set-TextValue $q1.Range("D2") "2.83"
Set-TextValue $q1.Range("E2") "79.91"
Set-TextValue $q1.Range("F2") "7.85"
Set-TextValue $q1.Range("G2") "0.2222"
$q1.Range("H2").Value = 3

$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "011588"
Set-TextValue $q1.Range("C3") "前海开源成份精选混合"
Set-TextValue $q1.Range("D3") "0.96"
Set-TextValue $q1.Range("E3") "81.61"
Set-TextValue $q1.Range("F3") "6.51"
Set-TextValue $q1.Range("G3") "0.0625"
$q1.Range("H3").Value = 5

# Apply the same header/index-column style used by the other quarterly
# sheets (bold, centred, bordered) by copying it over from "2021-Q4".
$fmtSrc = $wb.Worksheets.Item("2021-Q4")
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Insert a brand-new "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add($null, $q1)
$newTotal.Name = "总计"
$newTotal.Outline.SummaryBelow = $true
$newTotal.Outline.SummaryRight = $true
$newTotal.StandardWidth = 8.43

Set-TextValue $newTotal.Range("B1") "日期"
Set-TextValue $newTotal.Range("C1") "持有数量(只)"
Set-TextValue $newTotal.Range("D1") "持有市值(亿元)"

$newTotal.Range("A2").Value = 0
Set-TextValue $newTotal.Range("B2") "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 0.28

$newTotal.Range("A3").Value = 1
Set-TextValue $newTotal.Range("B3") "2021-Q4"
$newTotal.Range("C3").Value = 9
$newTotal.Range("D3").Value = 0.98

$newTotal.Range("A4").Value = 2
Set-TextValue $newTotal.Range("B4") "2021-Q3"
$newTotal.Range("C4").Value = 13
$newTotal.Range("D4").Value = 3.13

$newTotal.Range("A5").Value = 3
Set-TextValue $newTotal.Range("B5") "2021-Q2"
$newTotal.Range("C5").Value = 10
$newTotal.Range("D5").Value = 1.25

$newTotal.Range("A6").Value = 4
Set-TextValue $newTotal.Range("B6") "2021-Q1"
$newTotal.Range("C6").Value = 11
$newTotal.Range("D6").Value = 1.36

$newTotal.Range("A7").Value = 5
Set-TextValue $newTotal.Range("B7") "2020-Q4"
$newTotal.Range("C7").Value = 2
$newTotal.Range("D7").Value = 0.09

# Apply the same header/index-column style used before on the "总计" sheet.
$fmtSrc.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A3").Copy()
$newTotal.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Worksheets.Item("2020-Q4").Select()
